$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Preserve current formatting of rows 3-6 (cols A:B) onto rows 7-10 before
#    we touch anything, so the new "answer" rows get the same date / text
#    cell styles (s=2 / s=3) that rows 3-6 currently use.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Clear out the old submissions from rows 3-6 and reset them to the plain
#    (unused-row) style, matching cells C3:H6.
# ---------------------------------------------------------------------------
$ws.Range("A3:B6").ClearContents()
$ws.Range("C3").Copy()
$ws.Range("A3:B6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Write the new submissions into rows 7-10.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 45285.016736111109
$ws.Range("B7").Value = "הקבוצה של: טון, לא נתמך עי גוגל, המפקד"

$ws.Range("A8").Value = 45285.016840277778
$ws.Range("B8").Value = "הקבוצה של: דור, אלכס"

$ws.Range("A9").Value = 45285.016932870371
$ws.Range("B9").Value = "הקבוצה של: איי, הקשבי"

$ws.Range("A10").Value = 45285.017025462963
$ws.Range("B10").Value = "הקבוצה של: אריה, עמרי"

# ---------------------------------------------------------------------------
# 4) Append four new empty (plain-style) rows at the bottom of the sheet,
#    growing the used range from H106 to H110.
# ---------------------------------------------------------------------------
$ws.Range("A103").Copy()
$ws.Range("A107:H110").PasteSpecial(-4122)
$ws.Range("A107:H110").RowHeight = 15.75

# ---------------------------------------------------------------------------
# 5) Move the active selection to D7 (was E12).
# ---------------------------------------------------------------------------
$ws.Range("D7").Select() | Out-Null

$excel.CutCopyMode = 0
